$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = '26.792.37'
$ws.Range("E2").Value = '  +4.04%  '
$ws.Range("D3").Value = '1.876.20'
$ws.Range("E3").Value = '  +3.50%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '278.40'
$ws.Range("E5").Value = '  +0.88%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  +0.15%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5252'
$ws.Range("E7").Value = '  +3.77%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3453'
$ws.Range("E8").Value = '  -1.58%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.14'
$ws.Range("E9").Value = '  +0.23%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06943'
$ws.Range("E10").Value = '  +4.09%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.07'
$ws.Range("E11").Value = '  +0.49%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.8043'
$ws.Range("E12").Value = '  -3.18%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07871'
$ws.Range("E13").Value = '  -0.42%  '
$ws.Range("D14").Value = '1.891.56'
$ws.Range("E14").Value = '  +4.39%  '
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.152'
$ws.Range("E15").Value = '  +1.53%  '
$ws.Range("B16").Value = 'Litecoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '89.66'
$ws.Range("E16").Value = '  +2.55%  '
$ws.Range("B17").Value = 'BinanceUSD'
$ws.Range("C17").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.003'
$ws.Range("E17").Value = '  +0.21%  '
$ws.Range("B18").Value = 'Avalanche'
$ws.Range("C18").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.52'
$ws.Range("E18").Value = '  +3.95%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000008066'
$ws.Range("E19").Value = '  +0.45%  '
$ws.Range("E20").Value = '  +0.11%  '
$ws.Range("D21").Value = '26.861.03'
$ws.Range("E21").Value = '  +4.13%  '
$ws.Range("D22").Value = '2.128.01'
$ws.Range("E22").Value = '  +4.53%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.748'
$ws.Range("E23").Value = '  +0.67%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.02'
$ws.Range("E24").Value = '  +0.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '6.184'
$ws.Range("E25").Value = '  +2.42%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.337'
$ws.Range("E26").Value = '  +7.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '146.69'
$ws.Range("E27").Value = '  +3.64%  '
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.671'
$ws.Range("E28").Value = '  +0.27%  '
$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '17.39'
$ws.Range("E29").Value = '  +2.16%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '113.82'
$ws.Range("E30").Value = '  +4.03%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.354'
$ws.Range("E31").Value = '  +0.87%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.322'
$ws.Range("E32").Value = '  +2.29%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.08912'
$ws.Range("E33").Value = '  +1.39%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04955'
$ws.Range("E34").Value = '  +1.95%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.172'
$ws.Range("E35").Value = '  +3.25%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7343'
$ws.Range("E36").Value = '  +1.32%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.907'
$ws.Range("E37").Value = '  +1.17%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.284'
$ws.Range("E38").Value = '  +4.53%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.388'
$ws.Range("E39").Value = '  +5.82%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01849'
$ws.Range("E40").Value = '  +0.60%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5129'
$ws.Range("E41").Value = '  -1.14%  '
$ws.Range("E42").Value = '  +1.19%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '116.17'
$ws.Range("E43").Value = '  +2.72%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '6.208'
$ws.Range("E44").Value = '  +1.23%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.064'
$ws.Range("E45").Value = '  +0.20%  '
$ws.Range("E46").Value = '  +0.12%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4507'
$ws.Range("E47").Value = '  -0.91%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1342'
$ws.Range("E48").Value = '  -1.09%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.397'
$ws.Range("E49").Value = '  +1.13%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '36.53'
$ws.Range("E50").Value = '  +0.93%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05945'
$ws.Range("E51").Value = '  +2.05%  '

